# Update the public EPEX Spot prices workbook:
#  - "Prix Spot": append a new day column (BF) with header "10-aug" and
#    one value per hourly row (2-25).
#  - "Gaz" and "CO2": append a new day row (55) for 2025-08-08.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": new column BF ("10-aug")
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell - copy the formatting of the previous header (BE1) so the
# new header cell keeps the same bold / centered / bordered style.
$wsPrix.Range("BF1").Value = "10-aug"
$wsPrix.Range("BE1").Copy()
$wsPrix.Range("BF1").PasteSpecial(-4122)   # xlPasteFormats

$prixSpotValues = @{
    2  = 77.88
    3  = 67.28
    4  = 58.05
    5  = 49.32
    6  = 46.22
    7  = 46.01
    8  = 42.41
    9  = 49.65
    10 = 26.87
    11 = -0.01
    12 = -2.26
    13 = -15.6
    14 = -14.05
    15 = -36.19
    16 = -50.29
    17 = -21
    18 = -1.16
    19 = 1.72
    20 = 26.5
    21 = 80.5
    22 = 98.01000000000001
    23 = 96.14
    24 = 90.59
    25 = 75.56999999999999
}

foreach ($row in $prixSpotValues.Keys) {
    $wsPrix.Cells.Item($row, 58).Value = $prixSpotValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": new row 55 (2025-08-08 / 31.2)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to text first so Excel does not auto-convert the
# "2025-08-08" string into a date serial number (matches the existing
# rows, which store the date as plain text).
$wsGaz.Range("A55").NumberFormat = "@"
$wsGaz.Range("A55").Value = "2025-08-08"
$wsGaz.Range("A55").Style = "Normal"
$wsGaz.Range("B55").Value = 31.2

# ---------------------------------------------------------------------
# Sheet "CO2": new row 55 (2025-08-08 / 71.75)
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A55").NumberFormat = "@"
$wsCo2.Range("A55").Value = "2025-08-08"
$wsCo2.Range("A55").Style = "Normal"
$wsCo2.Range("B55").Value = 71.75
